$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '58.847.91'
Set-TextValue 'E2' '  +2.34%  '

Set-TextValue 'D3' '2.545.29'
Set-TextValue 'E3' '  +4.39%  '

Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.12%  '

Set-TextValue 'D5' '527.79'
Set-TextValue 'E5' '  +2.74%  '

Set-TextValue 'D6' '134.64'
Set-TextValue 'E6' '  +2.18%  '

Set-TextValue 'E7' '  -0.53%  '

Set-TextValue 'D8' '0.567'
Set-TextValue 'E8' '  +2.52%  '

Set-TextValue 'D9' '2.543.80'
Set-TextValue 'E9' '  +4.35%  '

Set-TextValue 'E10' '  +2.91%  '

Set-TextValue 'E11' '  -0.87%  '

Set-TextValue 'D12' '5.21'
Set-TextValue 'E12' '  -0.36%  '

Set-TextValue 'E13' '  +1.33%  '

Set-TextValue 'D14' '2.993.35'
Set-TextValue 'E14' '  +3.77%  '

Set-TextValue 'D15' '59.089.48'
Set-TextValue 'E15' '  +2.64%  '

Set-TextValue 'D16' '22.46'
Set-TextValue 'E16' '  +4.04%  '

Set-TextValue 'D17' '0.0000137'
Set-TextValue 'E17' '  +3.05%  '

Set-TextValue 'D18' '2.538.49'
Set-TextValue 'E18' '  +3.44%  '

Set-TextValue 'E19' '  +3.36%  '

Set-TextValue 'D20' '324.40'
Set-TextValue 'E20' '  +2.83%  '

Set-TextValue 'E21' '  +2.64%  '

Set-TextValue 'D22' '6.18'
Set-TextValue 'E22' '  +9.72%  '

Set-TextValue 'E23' '  +0.11%  '

Set-TextValue 'D24' '65.38'
Set-TextValue 'E24' '  +2.13%  '

Set-TextValue 'E25' '  +2.09%  '

Set-TextValue 'E26' '  -0.40%  '

Set-TextValue 'E27' '  +0.86%  '

Set-TextValue 'D28' '7.49'
Set-TextValue 'E28' '  +3.61%  '

Set-TextValue 'D29' '0.0₃0758'
Set-TextValue 'E29' '  +3.91%  '

Set-TextValue 'B30' 'Fetch.AI'
Set-TextValue 'C30' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D30' '1.22'
Set-TextValue 'E30' '  +5.77%  '

Set-TextValue 'B31' 'PancakeSwap'
Set-TextValue 'C31' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D31' '1.74'
Set-TextValue 'E31' '  +4.29%  '

Set-TextValue 'D32' '168.75'
Set-TextValue 'E32' '  -0.90%  '

Set-TextValue 'E33' '  +2.59%  '

Set-TextValue 'D34' '0.999'
Set-TextValue 'E34' '  -0.04%  '

Set-TextValue 'D35' '0.999'
Set-TextValue 'E35' '  -0.12%  '

Set-TextValue 'D36' '18.32'
Set-TextValue 'E36' '  +2.98%  '

Set-TextValue 'D37' '1.28'
Set-TextValue 'E37' '  -1.37%  '

Set-TextValue 'E38' '  +2.44%  '

Set-TextValue 'E39' '  +5.04%  '

Set-TextValue 'D40' '36.81'
Set-TextValue 'E40' '  +1.17%  '

Set-TextValue 'D41' '0.786'
Set-TextValue 'E41' '  +0.85%  '

Set-TextValue 'D42' '281.68'
Set-TextValue 'E42' '  +5.07%  '

Set-TextValue 'E43' '  +3.58%  '

Set-TextValue 'D44' '133.89'
Set-TextValue 'E44' '  +9.50%  '

Set-TextValue 'D45' '5.10'
Set-TextValue 'E45' '  +3.06%  '

Set-TextValue 'D46' '0.603'
Set-TextValue 'E46' '  +3.59%  '

Set-TextValue 'E47' '  +2.42%  '

Set-TextValue 'D48' '0.0508'
Set-TextValue 'E48' '  +5.20%  '

Set-TextValue 'E49' '  +3.63%  '

Set-TextValue 'D50' '0.0218'
Set-TextValue 'E50' '  +3.60%  '

Set-TextValue 'D51' '17.21'
Set-TextValue 'E51' '  +3.37%  '
